$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 83  # F2: 82 -> 83
$ws.Cells.Item(4, 6).Value = 2121  # F4: 2119 -> 2121
$ws.Cells.Item(8, 6).Value = 2102  # F8: 2101 -> 2102
$ws.Cells.Item(10, 6).Value = 10921  # F10: 10914 -> 10921
$ws.Cells.Item(12, 6).Value = 163  # F12: 162 -> 163
$ws.Cells.Item(14, 6).Value = 209  # F14: 208 -> 209
$ws.Cells.Item(15, 6).Value = 10717  # F15: 10712 -> 10717
$ws.Cells.Item(17, 6).Value = 1125  # F17: 1124 -> 1125
$ws.Cells.Item(18, 6).Value = 5  # F18: 4 -> 5
$ws.Cells.Item(19, 6).Value = 747  # F19: 745 -> 747
$ws.Cells.Item(20, 6).Value = 5361  # F20: 5359 -> 5361

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 27  # F2: 26 -> 27
$ws.Cells.Item(3, 6).Value = 561  # F3: 560 -> 561

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 83  # F2: 82 -> 83
$ws.Cells.Item(4, 6).Value = 2121  # F4: 2119 -> 2121
$ws.Cells.Item(7, 6).Value = 27  # F7: 26 -> 27
$ws.Cells.Item(9, 6).Value = 2102  # F9: 2101 -> 2102
$ws.Cells.Item(10, 6).Value = 561  # F10: 560 -> 561
$ws.Cells.Item(13, 6).Value = 10921  # F13: 10914 -> 10921
$ws.Cells.Item(15, 6).Value = 163  # F15: 162 -> 163
$ws.Cells.Item(17, 6).Value = 209  # F17: 208 -> 209
$ws.Cells.Item(18, 6).Value = 10717  # F18: 10712 -> 10717
$ws.Cells.Item(20, 6).Value = 1125  # F20: 1124 -> 1125
$ws.Cells.Item(21, 6).Value = 5  # F21: 4 -> 5
$ws.Cells.Item(22, 6).Value = 747  # F22: 745 -> 747
$ws.Cells.Item(23, 6).Value = 5361  # F23: 5359 -> 5361
